$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6: "Change in inventories"
$ws.Range("B6").Value = -275600000.0
$ws.Range("C6").Value = -145700000.0
$ws.Range("D6").Value = -96800000.0
$ws.Range("E6").Value = -71000000.0
$ws.Range("F6").Value = -123400000.0

# Row 7: "Change in payables and accrued liability"
$ws.Range("B7").Value = 367300000.0
$ws.Range("C7").Value = 272300000.0
$ws.Range("D7").Value = 133400000.0
$ws.Range("E7").Value = 100600000.0
$ws.Range("F7").Value = 93200000.0
